$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.162.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.35%  "

$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5179"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.83%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4013"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08442"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.119"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.427"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.906.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.333"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001110"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06677"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "

$ws.Range("E20").Value = "  +2.74%  "

$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.945"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.173.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.214"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.120.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.388"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.095"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.058"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.727"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02489"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06562"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2211"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.256"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.217"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.757"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6501"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.237"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6106"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.707"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.056"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.235"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.166"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.18%  "
